# The sheet contains one data row per market report. A new daily/weekly
# report needs to be inserted as the new row 62 (pushing the existing
# rows 62-159 down to 63-160, so the sheet grows from A1:R159 to A1:R160).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 62, shifting rows 62..159 down to 63..160.
$ws.Rows.Item(62).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row 62 with the new report's data.
$ws.Cells.Item(62, 1).Value  = 4
$ws.Cells.Item(62, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(62, 3).Value  = "Los Lagos"
$ws.Cells.Item(62, 4).Value  = 44571
$ws.Cells.Item(62, 5).Value  = 10
$ws.Cells.Item(62, 6).Value  = 100112039
$ws.Cells.Item(62, 7).Value  = "Ciboulette"
$ws.Cells.Item(62, 8).Value  = "Sin especificar"
$ws.Cells.Item(62, 9).Value  = "Primera"
$ws.Cells.Item(62, 10).Value = 80
$ws.Cells.Item(62, 11).Value = 3000
$ws.Cells.Item(62, 12).Value = 3000
$ws.Cells.Item(62, 13).Value = 3000
$ws.Cells.Item(62, 14).Value = "$/docena de atados"
$ws.Cells.Item(62, 15).Value = "Región Metropolitana"
$ws.Cells.Item(62, 16).Value = 1000
$ws.Cells.Item(62, 17).Value = 3
$ws.Cells.Item(62, 18).Value = "Hortaliza"
